$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the missing "Fecha de Inicio" / "Fecha de Fin" dates for rows 16-20 ---
# Copy the existing date format (used on column H) into I16:J20 first, so the new
# values pick up the same number format / style as the rest of the date columns,
# then write the actual date serials.
$ws.Range("H16").Copy()
$ws.Range("I16:J20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I16").Value = 40886
$ws.Range("J16").Value = 40886
$ws.Range("I17").Value = 40886
$ws.Range("J17").Value = 40886
$ws.Range("I18").Value = 40884
$ws.Range("J18").Value = 40886
$ws.Range("I19").Value = 40884
$ws.Range("J19").Value = 40886
$ws.Range("I20").Value = 40884
$ws.Range("J20").Value = 40886

# --- Remove the underline formatting from G14 ("Jorge Alcantara") ---
$ws.Range("G14").Font.Underline = $false

# --- Move the active selection to G14 ---
$ws.Range("G14").Select()
